$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 to the new numeric value.
$ws.Cells.Item(2, 2).Value = 1684.22140908781

# The "max" column (C) is dropped; "prediction" (old D) shifts into C,
# and "rejection-f" (old E) shifts into D. Copy header text first.
$ws.Cells.Item(1, 3).Value = $ws.Cells.Item(1, 4).Text
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(2, 4).Text

# Now remove the old column D's duplicate data by deleting column D
# (which shifts old E into D), leaving C already holding the old D
# header/value copied above.
$ws.Columns.Item(4).Delete()
